# Remove the review comment that was left on slide 1
# ("Button not visible. Not aligned with anything.").
#
# Deleting the Comment object removes ppt/comments/comment1.xml (and the
# now-dangling Content_Types override / slide relationship that pointed to
# it) while leaving ppt/commentAuthors.xml and everything else untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Comments.Count; $i -ge 1; $i--) {
    $s.Comments.Item($i).Delete()
}
